$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: AD1 = "Wins", AE1 = "Losses", AF1 = "Ties" ---
# Copy the formatting (bold font, border, center/top alignment) from an
# existing header cell (A1) so the new header cells match the rest of row 1,
# then set the text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows 2-48: team record columns (same constant values for every row) ---
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 75   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 87   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
